$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 19.81693648186462
$ws.Range("E2").Value = 19.35400009155273
$ws.Range("F2").Value = 20.45448944734193
$ws.Range("G2").Value = 19.09410608640153
$ws.Range("H2").Value = 539679667
$ws.Range("I2").Value = "MCHP"

$ws.Range("D3").Value = 19.59466919744172
$ws.Range("E3").Value = 17.52475166320801
$ws.Range("F3").Value = 19.6192126687538
$ws.Range("G3").Value = 16.73523775243184
$ws.Range("H3").Value = 539679667
$ws.Range("I3").Value = "MCHP"

$ws.Range("D4").Value = 17.76883405873871
$ws.Range("E4").Value = 19.91313552856445
$ws.Range("F4").Value = 20.45333346668414
$ws.Range("G4").Value = 17.05131721081299
$ws.Range("H4").Value = 539679667
$ws.Range("I4").Value = "MCHP"

$ws.Range("D5").Value = 18.97645025745353
$ws.Range("E5").Value = 18.62318801879883
$ws.Range("F5").Value = 19.2465930884726
$ws.Range("G5").Value = 16.47451772091467
$ws.Range("H5").Value = 539679667
$ws.Range("I5").Value = "MCHP"

$ws.Range("D6").Value = 20.08054212786136
$ws.Range("E6").Value = 20.36555099487305
$ws.Range("F6").Value = 21.24153380107096
$ws.Range("G6").Value = 19.91708148919265
$ws.Range("H6").Value = 539679667
$ws.Range("I6").Value = "MCHP"

$ws.Range("D7").Value = 21.36294229748012
$ws.Range("E7").Value = 23.49543571472168
$ws.Range("F7").Value = 24.80449224404243
$ws.Range("G7").Value = 20.6788556479555
$ws.Range("H7").Value = 539679667
$ws.Range("I7").Value = "MCHP"

$ws.Range("D8").Value = 26.36215604546758
$ws.Range("E8").Value = 25.72072982788086
$ws.Range("F8").Value = 26.59153869969209
$ws.Range("G8").Value = 24.43363187638488
$ws.Range("H8").Value = 539679667
$ws.Range("I8").Value = "MCHP"

$ws.Range("D9").Value = 27.55439599915982
$ws.Range("E9").Value = 28.76745414733887
$ws.Range("F9").Value = 29.5106658755721
$ws.Range("G9").Value = 26.57198718209643
$ws.Range("H9").Value = 539679667
$ws.Range("I9").Value = "MCHP"

$ws.Range("D10").Value = 31.74164841277895
$ws.Range("E10").Value = 32.44568252563477
$ws.Range("F10").Value = 33.19264490874432
$ws.Range("G10").Value = 30.65125331720539
$ws.Range("H10").Value = 539679667
$ws.Range("I10").Value = "MCHP"

$ws.Range("D11").Value = 33.5215896771575
$ws.Range("E11").Value = 34.51779174804688
$ws.Range("F11").Value = 35.95818775893925
$ws.Range("G11").Value = 32.68926399192877
$ws.Range("H11").Value = 539679667
$ws.Range("I11").Value = "MCHP"

$ws.Range("D12").Value = 39.14103643925164
$ws.Range("E12").Value = 41.05975723266602
$ws.Range("F12").Value = 41.44523357971038
$ws.Range("G12").Value = 38.9894418698669
$ws.Range("H12").Value = 539679667
$ws.Range("I12").Value = "MCHP"

$ws.Range("D13").Value = 38.51890408079754
$ws.Range("E13").Value = 41.40629959106445
$ws.Range("F13").Value = 43.12395087153036
$ws.Range("G13").Value = 38.37975260985574
$ws.Range("H13").Value = 539679667
$ws.Range("I13").Value = "MCHP"

$ws.Range("D14").Value = 39.73983479199187
$ws.Range("E14").Value = 36.5384635925293
$ws.Range("F14").Value = 39.96257796281496
$ws.Range("G14").Value = 35.49899768344432
$ws.Range("H14").Value = 539679667
$ws.Range("I14").Value = "MCHP"

$ws.Range("D15").Value = 39.55775010474248
$ws.Range("E15").Value = 40.96520233154297
$ws.Range("F15").Value = 42.49104029756656
$ws.Range("G15").Value = 39.110520293505
$ws.Range("H15").Value = 539679667
$ws.Range("I15").Value = "MCHP"

$ws.Range("D16").Value = 34.92292524220144
$ws.Range("E16").Value = 28.96885108947754
$ws.Range("F16").Value = 35.09907853862872
$ws.Range("G16").Value = 26.73167096175434
$ws.Range("H16").Value = 539679667
$ws.Range("I16").Value = "MCHP"

$ws.Range("D17").Value = 31.0998452125534
$ws.Range("E17").Value = 35.56986999511719
$ws.Range("F17").Value = 36.17177422866585
$ws.Range("G17").Value = 29.28085913544141
$ws.Range("H17").Value = 539679667
$ws.Range("I17").Value = "MCHP"

$ws.Range("D18").Value = 37.38885107756261
$ws.Range("E18").Value = 44.38759613037109
$ws.Range("F18").Value = 44.56089811760286
$ws.Range("G18").Value = 37.37552067704233
$ws.Range("H18").Value = 539679667
$ws.Range("I18").Value = "MCHP"

$ws.Range("D19").Value = 40.28546342551344
$ws.Range("E19").Value = 42.14218139648438
$ws.Range("F19").Value = 44.65946461999538
$ws.Range("G19").Value = 38.88845817131561
$ws.Range("H19").Value = 539679667
$ws.Range("I19").Value = "MCHP"

$ws.Range("D20").Value = 42.5158204809004
$ws.Range("E20").Value = 42.25589370727539
$ws.Range("F20").Value = 43.85578262629921
$ws.Range("G20").Value = 39.92552104234334

$ws.Range("D21").Value = 47.70615065146252
$ws.Range("E21").Value = 43.85924530029297
$ws.Range("F21").Value = 50.60370540202903
$ws.Range("G21").Value = 43.61178204698729
$ws.Range("H21").Value = 539679667
$ws.Range("I21").Value = "MCHP"

$ws.Range("D22").Value = 29.38870009393426
$ws.Range("E22").Value = 39.60477447509766
$ws.Range("F22").Value = 41.85745579806255
$ws.Range("G22").Value = 27.81769141081628
$ws.Range("H22").Value = 539679667
$ws.Range("I22").Value = "MCHP"

$ws.Range("D23").Value = 47.65880383559927
$ws.Range("E23").Value = 46.11747741699219
$ws.Range("F23").Value = 49.92093357164411
$ws.Range("G23").Value = 45.02947965312577
$ws.Range("H23").Value = 539679667
$ws.Range("I23").Value = "MCHP"

$ws.Range("D24").Value = 48.25370847608574
$ws.Range("E24").Value = 47.798828125
$ws.Range("F24").Value = 53.02085261200573
$ws.Range("G24").Value = 46.55700323953281
$ws.Range("H24").Value = 539679667
$ws.Range("I24").Value = "MCHP"

$ws.Range("D25").Value = 63.3929779055899
$ws.Range("E25").Value = 62.0972900390625
$ws.Range("F25").Value = 70.87969271257363
$ws.Range("G25").Value = 60.46855020172884
$ws.Range("H25").Value = 539679667
$ws.Range("I25").Value = "MCHP"

$ws.Range("D26").Value = 72.03548133110579
$ws.Range("E26").Value = 68.73785400390625
$ws.Range("F26").Value = 75.22790896235873
$ws.Range("G26").Value = 67.64474443937054
$ws.Range("H26").Value = 539679667
$ws.Range("I26").Value = "MCHP"

$ws.Range("D27").Value = 68.72971693859773
$ws.Range("E27").Value = 65.64295196533203
$ws.Range("F27").Value = 69.20213283225658
$ws.Range("G27").Value = 59.19423916586787
$ws.Range("H27").Value = 539679667
$ws.Range("I27").Value = "MCHP"

$ws.Range("D28").Value = 71.12469442965109
$ws.Range("E28").Value = 68.17112731933594
$ws.Range("F28").Value = 71.16609780468053
$ws.Range("G28").Value = 64.2790556287501
$ws.Range("H28").Value = 539679667
$ws.Range("I28").Value = "MCHP"

$ws.Range("D29").Value = 80.75434739721622
$ws.Range("E29").Value = 71.49048614501953
$ws.Range("F29").Value = 82.09226076423309
$ws.Range("G29").Value = 63.45379813105307
$ws.Range("H29").Value = 539679667
$ws.Range("I29").Value = "MCHP"

$ws.Range("D30").Value = 70.0694899585865
$ws.Range("E30").Value = 60.37439727783203
$ws.Range("F30").Value = 70.36580575012617
$ws.Range("G30").Value = 58.65206306414378
$ws.Range("H30").Value = 539679667
$ws.Range("I30").Value = "MCHP"

$ws.Range("D31").Value = 53.18845437543131
$ws.Range("E31").Value = 64.03071594238281
$ws.Range("F31").Value = 64.44915373663841
$ws.Range("G31").Value = 50.51973400462681
$ws.Range("H31").Value = 539679667
$ws.Range("I31").Value = "MCHP"

$ws.Range("D32").Value = 58.09098928328942
$ws.Range("E32").Value = 57.6521110534668
$ws.Range("F32").Value = 63.86180324860094
$ws.Range("G32").Value = 51.14360322232878
$ws.Range("H32").Value = 539679667
$ws.Range("I32").Value = "MCHP"

$ws.Range("D33").Value = 66.86548830200324
$ws.Range("E33").Value = 72.80262756347656
$ws.Range("F33").Value = 73.72180124862521
$ws.Range("G33").Value = 64.1642339860382
$ws.Range("H33").Value = 539679667
$ws.Range("I33").Value = "MCHP"

$ws.Range("D34").Value = 78.69745525351178
$ws.Range("E34").Value = 68.75077056884766
$ws.Range("F34").Value = 79.43215204333096
$ws.Range("G34").Value = 66.21700948425107
$ws.Range("H34").Value = 539679667
$ws.Range("I34").Value = "MCHP"

$ws.Range("D35").Value = 85.86158439712732
$ws.Range("E35").Value = 88.9189453125
$ws.Range("F35").Value = 89.25970402819497
$ws.Range("G35").Value = 81.31814689053806
$ws.Range("H35").Value = 539679667
$ws.Range("I35").Value = "MCHP"

$ws.Range("D36").Value = 74.48051994898105
$ws.Range("E36").Value = 67.82979583740234
$ws.Range("F36").Value = 77.7059788287086
$ws.Range("G36").Value = 65.85075112243541
$ws.Range("H36").Value = 539679667
$ws.Range("I36").Value = "MCHP"

$ws.Range("D37").Value = 85.25523184900275
$ws.Range("E37").Value = 81.46781158447266
$ws.Range("F37").Value = 88.66964910014818
$ws.Range("G37").Value = 78.15860105986154
$ws.Range("H37").Value = 539679667
$ws.Range("I37").Value = "MCHP"

$ws.Range("D38").Value = 86.29424630654887
$ws.Range("E38").Value = 88.4482421875
$ws.Range("F38").Value = 90.99648412484908
$ws.Range("G38").Value = 77.79367978671276
$ws.Range("H38").Value = 539679667
$ws.Range("I38").Value = "MCHP"

$ws.Range("D39").Value = 88.40091618046965
$ws.Range("E39").Value = 85.77304077148438
$ws.Range("F39").Value = 93.69531309221804
$ws.Range("G39").Value = 80.13084375771258
$ws.Range("H39").Value = 539679667
$ws.Range("I39").Value = "MCHP"

$ws.Range("D40").Value = 77.70676833571524
$ws.Range("E40").Value = 71.27573394775391
$ws.Range("F40").Value = 78.0370672508877
$ws.Range("G40").Value = 70.59571025727955
$ws.Range("H40").Value = 539679667
$ws.Range("I40").Value = "MCHP"

$ws.Range("D41").Value = 56.76190899844286
$ws.Range("E41").Value = 53.11341857910156
$ws.Range("F41").Value = 58.5519210540414
$ws.Range("G41").Value = 52.80041161823252
$ws.Range("H41").Value = 539679667
$ws.Range("I41").Value = "MCHP"

$ws.Range("D42").Value = 47.71742775524408
$ws.Range("E42").Value = 45.41138458251953
$ws.Range("F42").Value = 48.75219534425183
$ws.Range("G42").Value = 33.63477740275318
$ws.Range("H42").Value = 539679667
$ws.Range("I42").Value = "MCHP"

$ws.Range("D43").Value = 68.98187118333068
$ws.Range("E43").Value = 67.12474060058594
$ws.Range("F43").Value = 76.66859077208977
$ws.Range("G43").Value = 66.27066738128168
$ws.Range("H43").Value = 539679667
$ws.Range("I43").Value = "MCHP"
